$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 380. Excel shifts the
# existing rows 380-386 down to 382-388, preserving their content.
$ws.Rows.Item(380).Insert()
$ws.Rows.Item(380).Insert()

# --- New row 380: Mandarina Murcott, Primera, Región de O'Higgins ---
$ws.Cells.Item(380, 1).Value = 5
$ws.Cells.Item(380, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(380, 3).Value = "Maule"
$ws.Cells.Item(380, 4).Value = 44890
$ws.Cells.Item(380, 5).Value = 7
$ws.Cells.Item(380, 6).Value = "Fruta"
$ws.Cells.Item(380, 7).Value = 100102
$ws.Cells.Item(380, 8).Value = "Cítricos"
$ws.Cells.Item(380, 9).Value = 100102004
$ws.Cells.Item(380, 10).Value = "Mandarina"
$ws.Cells.Item(380, 11).Value = "Murcott"
$ws.Cells.Item(380, 12).Value = "Primera"
$ws.Cells.Item(380, 13).Value = 180
$ws.Cells.Item(380, 14).Value = 6000
$ws.Cells.Item(380, 15).Value = 6000
$ws.Cells.Item(380, 16).Value = 6000
$ws.Cells.Item(380, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(380, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(380, 19).Value = 333
$ws.Cells.Item(380, 20).Value = 18

# --- New row 381: Mandarina Murcott, Segunda, Región de O'Higgins ---
$ws.Cells.Item(381, 1).Value = 5
$ws.Cells.Item(381, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(381, 3).Value = "Maule"
$ws.Cells.Item(381, 4).Value = 44890
$ws.Cells.Item(381, 5).Value = 7
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100102
$ws.Cells.Item(381, 8).Value = "Cítricos"
$ws.Cells.Item(381, 9).Value = 100102004
$ws.Cells.Item(381, 10).Value = "Mandarina"
$ws.Cells.Item(381, 11).Value = "Murcott"
$ws.Cells.Item(381, 12).Value = "Segunda"
$ws.Cells.Item(381, 13).Value = 150
$ws.Cells.Item(381, 14).Value = 5000
$ws.Cells.Item(381, 15).Value = 5000
$ws.Cells.Item(381, 16).Value = 5000
$ws.Cells.Item(381, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(381, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(381, 19).Value = 278
$ws.Cells.Item(381, 20).Value = 18

# Make sure the date cells keep the same date-time number format as the
# rest of column D (style index used throughout column D).
$ws.Range("D380:D381").NumberFormat = $ws.Range("D382").NumberFormat
